# Update the "Expertise" (column G) ratings on the marker sheet and move
# the current selection, matching the author's manual review pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 Yetta Cisneros")

# Select the row the reviewer was working on (A10:K10), matching the new
# selection state left in the worksheet after editing.
[void]$ws.Range("A10:K10").Select()

# Clear / change / fill in the expertise ratings for the affected rows.
$ws.Range("G10").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("G12").Value = "M"
$ws.Range("G13").Value = "H"
$ws.Range("G17").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("G19").Value = "L"
$ws.Range("G21").Value = "L"
$ws.Range("G22").Value = "L"
$ws.Range("G23").Value = ""
